$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.345.80"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.79%  "

$ws.Range("D3").Value = "'3.708.15"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.20%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'596.67"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.45%  "

$ws.Range("D6").Value = "'166.29"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.37%  "

$ws.Range("D7").Value = "'3.712.07"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.16%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("E10").Value = "  -3.01%  "

$ws.Range("D11").Value = "'6.18"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.62%  "

$ws.Range("E12").Value = "  -4.33%  "

$ws.Range("D13").Value = "'37.86"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.79%  "

$ws.Range("E14").Value = "  -5.26%  "

$ws.Range("D15").Value = "'4.322.82"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.73%  "

$ws.Range("D16").Value = "'3.701.68"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.60%  "

$ws.Range("D17").Value = "'67.422.12"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.82%  "

$ws.Range("D18").Value = "'17.56"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.04%  "

$ws.Range("D19").Value = "'7.21"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.96%  "

$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").Value = "'488.58"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.78%  "

$ws.Range("D22").Value = "'9.35"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.37%  "

$ws.Range("D23").Value = "'0.726"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").Value = "'85.43"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.44%  "

$ws.Range("E25").Value = "  -7.20%  "

$ws.Range("E26").Value = "  -4.77%  "

$ws.Range("D27").Value = "'12.21"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.73%  "

$ws.Range("D28").Value = "'10.14"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.30%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  -2.25%  "

$ws.Range("D31").Value = "'2.36"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.55%  "

$ws.Range("D32").Value = "'7.67"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.72%  "

$ws.Range("D33").Value = "'31.50"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.91%  "

$ws.Range("D34").Value = "'3.842.69"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.03%  "

$ws.Range("E35").Value = "  -5.33%  "

$ws.Range("D36").Value = "'3.646.85"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.12%  "

$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.52%  "

$ws.Range("D39").Value = "'5.75"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.61%  "

$ws.Range("E40").Value = "  -7.78%  "

$ws.Range("D41").Value = "'0.323"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.59%  "

$ws.Range("D42").Value = "'429.84"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.50%  "

$ws.Range("D43").Value = "'48.60"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("E44").Value = "  -6.51%  "

$ws.Range("E45").Value = "  -6.33%  "

$ws.Range("D46").Value = "'8.43"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "'40.62"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.55%  "

$ws.Range("D49").Value = "'141.42"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("D50").Value = "'2.752.26"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.57%  "

$ws.Range("E51").Value = "  -4.39%  "
